$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 4) with a TimeStamp header/value pair.
$ws.Cells.Item(4, 1).Value = "TimeStamp"
$ws.Cells.Item(4, 2).Value = "11/26/2020 2:42:15 AM"
$ws.Cells.Item(4, 3).Value = "11/26/2020 2:42:15 AM"

$ws.Range("C4").Select()
